$d = $word.ActiveDocument

# Locate the target paragraph
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Issue 1:*") {
        $para = $p
        break
    }
}
$pStart = $para.Range.Start

# Helper: force a run boundary at an absolute document position by
# splitting the paragraph there and immediately deleting the paragraph
# mark again (merges the paragraphs back but keeps two separate runs,
# without leaving any residual run formatting / rPr artifacts).
function Split-RunAt([int]$pos) {
    $ins = $d.Range($pos, $pos)
    $ins.InsertParagraphAfter()
    $mark = $d.Range($pos, $pos + 1)
    $mark.Delete()
}

# --- Step 1: perform all text content edits first --------------------

# Remove the trailing period.
$periodRange = $d.Range($pStart + 128, $pStart + 129)
$periodRange.Delete()

# Insert " (Potentially Fixed)" right after "Issue 1".
$insPoint = $d.Range($pStart + 7, $pStart + 7)
$insPoint.InsertAfter(" (Potentially Fixed)")

$full = $d.Content.Text
$idx = $full.IndexOf("Issue 1")
Write-Host "AFTER TEXT EDITS=[$($full.Substring($idx, 150))]"

# --- Step 2: now that text content is final, force the run boundaries
#             (rightmost first so earlier offsets remain valid) -------

Split-RunAt ($pStart + 117)
Split-RunAt ($pStart + 27)
Split-RunAt ($pStart + 7)

$full2 = $d.Content.Text
$idx2 = $full2.IndexOf("Issue 1")
Write-Host "RESULT=[$($full2.Substring($idx2, 160))]"
